# Auto-generated edit script: update crypto Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.720.45"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "2.247.23"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "116.08"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "301.53"
$ws.Range("E6").Value = "  +13.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  -2.32%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.625"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("E10").Value = "  -2.10%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.85"
$ws.Range("E12").Value = "  +4.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.17"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.900"
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("D17").Value = "2.587.69"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "2.271.14"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").Value = "42.857.01"
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.78"
$ws.Range("E20").Value = "  +12.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000109"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.08"
$ws.Range("E22").Value = "  +2.20%  "
$ws.Range("E23").Value = "  +25.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.38"
$ws.Range("E24").Value = "  -3.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "233.26"
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.50"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.23"
$ws.Range("E27").Value = "  +5.29%  "
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.53"
$ws.Range("E29").Value = "  -3.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.24"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("E31").Value = "  -4.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "175.83"
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.37"
$ws.Range("E33").Value = "  -1.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0913"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.60"
$ws.Range("E35").Value = "  +17.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.71"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.78"
$ws.Range("E38").Value = "  +1.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0377"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.60"
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.241"
$ws.Range("E42").Value = "  +1.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.91"
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.53"
$ws.Range("E44").Value = "  -5.50%  "
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.60"
$ws.Range("E47").Value = "  -6.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.37"
$ws.Range("E48").Value = "  +5.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "108.32"
$ws.Range("E49").Value = "  +7.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.66"
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("E51").Value = "  +7.98%  "
